$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.007.31"
$ws.Range("E2").Value = "  -3.80%  "

$ws.Range("D3").Value = "1.644.54"
$ws.Range("E3").Value = "  -5.59%  "

$ws.Range("D4").Value = "'0.9966"
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").Value = "'233.64"
$ws.Range("E5").Value = "  -5.13%  "

$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").Value = "'0.4797"
$ws.Range("E7").Value = "  -4.98%  "

$ws.Range("D8").Value = "'0.2602"
$ws.Range("E8").Value = "  -5.19%  "

$ws.Range("D9").Value = "'0.06107"
$ws.Range("E9").Value = "  -1.18%  "

$ws.Range("D10").Value = "'0.07075"
$ws.Range("E10").Value = "  -2.40%  "

$ws.Range("D11").Value = "1.638.92"
$ws.Range("E11").Value = "  -5.98%  "

$ws.Range("D12").Value = "'14.63"
$ws.Range("E12").Value = "  -3.28%  "

$ws.Range("D13").Value = "'0.6019"
$ws.Range("E13").Value = "  -7.61%  "

$ws.Range("D14").Value = "'4.385"

$ws.Range("D15").Value = "'73.89"
$ws.Range("E15").Value = "  -4.69%  "

$ws.Range("D16").Value = "'0.9987"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "'0.9979"
$ws.Range("E17").Value = "  -0.30%  "

$ws.Range("D18").Value = "25.008.81"
$ws.Range("E18").Value = "  -3.85%  "

$ws.Range("E19").Value = "  -4.02%  "

$ws.Range("E20").Value = "  -5.10%  "

$ws.Range("D21").Value = "1.848.98"
$ws.Range("E21").Value = "  -6.28%  "

$ws.Range("D22").Value = "'4.395"
$ws.Range("E22").Value = "  -1.38%  "

$ws.Range("D23").Value = "'8.614"
$ws.Range("E23").Value = "  -1.13%  "

$ws.Range("D24").Value = "'5.256"
$ws.Range("E24").Value = "  -2.09%  "

$ws.Range("D25").Value = "'133.60"
$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("D26").Value = "'14.90"
$ws.Range("E26").Value = "  -2.33%  "

$ws.Range("D27").Value = "'1.385"
$ws.Range("E27").Value = "  -7.98%  "

$ws.Range("D28").Value = "'104.02"
$ws.Range("E28").Value = "  -1.67%  "

$ws.Range("D29").Value = "'1.648"
$ws.Range("E29").Value = "  -7.37%  "

$ws.Range("D30").Value = "'3.876"

$ws.Range("D31").Value = "'0.07708"
$ws.Range("E31").Value = "  -5.82%  "

$ws.Range("D32").Value = "'3.557"
$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("D33").Value = "'0.9975"
$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("D34").Value = "'0.04288"
$ws.Range("E34").Value = "  -8.01%  "

$ws.Range("D35").Value = "'2.568"
$ws.Range("E35").Value = "  -3.47%  "

$ws.Range("D36").Value = "'0.9299"
$ws.Range("E36").Value = "  -6.54%  "

$ws.Range("D37").Value = "'0.5895"
$ws.Range("E37").Value = "  -3.14%  "

$ws.Range("D38").Value = "'2.568"
$ws.Range("E38").Value = "  -8.09%  "

$ws.Range("D39").Value = "'0.01521"
$ws.Range("E39").Value = "  -6.25%  "

$ws.Range("D40").Value = "'0.9976"
$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("D41").Value = "'0.8336"
$ws.Range("E41").Value = "  +9.10%  "

$ws.Range("B42").Value = "PaxosStandard"
$ws.Range("C42").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D42").Value = "'0.9977"
$ws.Range("E42").Value = "  -0.40%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'98.52"
$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.757"
$ws.Range("E44").Value = "  -8.99%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.3704"
$ws.Range("E45").Value = "  -5.32%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'4.693"
$ws.Range("E46").Value = "  -6.20%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1095"
$ws.Range("E47").Value = "  -5.82%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'6.101"
$ws.Range("E48").Value = "  -3.25%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05202"
$ws.Range("E49").Value = "  -2.06%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'29.18"
$ws.Range("E50").Value = "  -4.68%  "

$ws.Range("B51").Value = "TrueUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D51").Value = "'0.9981"
$ws.Range("E51").Value = "  -0.44%  "

Write-Host "Update complete"